$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shortRange = "May 25 - May 28"
$combinedRange = "May 25 - May 28, Aug 4 - Aug 6"

# Row 7 records a single survey window; every other data row records both
# windows used so far. Write E7 first so its (shorter) string is the first
# new shared string, then fill the rest top-to-bottom.
$ws.Range("E7").Value = $shortRange

for ($r = 2; $r -le 21; $r++) {
    if ($r -ne 7) {
        $ws.Cells.Item($r, 5).Value = $combinedRange
    }
}

# The column header changes from "Occurrence" to "Dates Used"
$ws.Range("E1").Value = "Dates Used"

# The saved selection moves to E1
$ws.Range("E1").Select()
